# Add the new "TankCollection" worksheet (mirrors schema2.yaml / examples).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Add()
$ws.Name = "TankCollection"
$ws.Range("A1").Value = "tanks"

# Match the page-margin conventions used by the existing sheets
# (0.75/0.75/1/1 in, 0.5/0.5 in header/footer -> points).
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Place the new sheet right after "Tank" (i.e. as the 3rd/last sheet).
$tankSheet = $wb.Worksheets.Item("Tank")
$ws.Move($null, $tankSheet)
